$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1836065573770492
$ws.Range("C2").Value = 0.5868852459016394
$ws.Range("J2").Value = 0.01311475409836066
$ws.Range("P2").Value = 0.1180327868852459
$ws.Range("S2").Value = 0.09836065573770492
$ws.Range("B3").Value = 0.01111111111111111
$ws.Range("C3").Value = 0.02222222222222222
$ws.Range("J3").Value = 0.02222222222222222
$ws.Range("P3").Value = 0.7388888888888889
$ws.Range("S3").Value = 0.2055555555555555
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7560975609756098
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.05785123966942149
$ws.Range("D6").Value = 0.008264462809917356
$ws.Range("F6").Value = 0.04132231404958678
$ws.Range("J6").Value = 0.268595041322314
$ws.Range("O6").Value = 0.01652892561983471
$ws.Range("Q6").Value = 0.1528925619834711
$ws.Range("R6").Value = 0.1074380165289256
$ws.Range("S6").Value = 0.3471074380165289
$ws.Range("B7").Value = 0.07741935483870968
$ws.Range("D7").Value = 0.01935483870967742
$ws.Range("E7").Value = 0.006451612903225806
$ws.Range("F7").Value = 0.08387096774193549
$ws.Range("J7").Value = 0.1548387096774194
$ws.Range("O7").Value = 0.02580645161290323
$ws.Range("Q7").Value = 0.1741935483870968
$ws.Range("R7").Value = 0.09032258064516129
$ws.Range("S7").Value = 0.367741935483871
$ws.Range("B8").Value = 0.09071274298056156
$ws.Range("D8").Value = 0.02159827213822894
$ws.Range("F8").Value = 0.06263498920086392
$ws.Range("J8").Value = 0.1252699784017278
$ws.Range("O8").Value = 0.01727861771058315
$ws.Range("Q8").Value = 0.16414686825054
$ws.Range("R8").Value = 0.142548596112311
$ws.Range("S8").Value = 0.3758099352051836
$ws.Range("B9").Value = 0.09777777777777778
$ws.Range("D9").Value = 0.03111111111111111
$ws.Range("F9").Value = 0.07555555555555556
$ws.Range("J9").Value = 0.1155555555555556
$ws.Range("O9").Value = 0.008888888888888889
$ws.Range("Q9").Value = 0.2044444444444445
$ws.Range("R9").Value = 0.09777777777777778
$ws.Range("S9").Value = 0.3688888888888889
$ws.Range("B10").Value = 0.1066098081023454
$ws.Range("D10").Value = 0.01350390902629709
$ws.Range("E10").Value = 0.002132196162046908
$ws.Range("F10").Value = 0.07604832977967306
$ws.Range("J10").Value = 0.1101634683724236
$ws.Range("O10").Value = 0.01137171286425018
$ws.Range("Q10").Value = 0.2217484008528785
$ws.Range("R10").Value = 0.09950248756218906
$ws.Range("S10").Value = 0.3589196872778962
$ws.Range("G11").Value = 0.148
$ws.Range("J11").Value = 0.116
$ws.Range("K11").Value = 0.232
$ws.Range("L11").Value = 0.484
$ws.Range("S11").Value = 0.02
$ws.Range("G12").Value = 0.71900826446281
$ws.Range("J12").Value = 0.2148760330578512
$ws.Range("L12").Value = 0.01652892561983471
$ws.Range("S12").Value = 0.04958677685950413
$ws.Range("G13").Value = 0.5909090909090909
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.07575757575757576
$ws.Range("F15").Value = 0.01310043668122271
$ws.Range("H15").Value = 0.1834061135371179
$ws.Range("I15").Value = 0.1091703056768559
$ws.Range("J15").Value = 0.4148471615720524
$ws.Range("K15").Value = 0.06986899563318777
$ws.Range("M15").Value = 0.02183406113537118
$ws.Range("O15").Value = 0.02620087336244541
$ws.Range("S15").Value = 0.1615720524017467
$ws.Range("F16").Value = 0.01015228426395939
$ws.Range("H16").Value = 0.1979695431472081
$ws.Range("I16").Value = 0.08629441624365482
$ws.Range("J16").Value = 0.3807106598984771
$ws.Range("K16").Value = 0.07614213197969544
$ws.Range("M16").Value = 0.03553299492385787
$ws.Range("N16").Value = 0.005076142131979695
$ws.Range("O16").Value = 0.05076142131979695
$ws.Range("S16").Value = 0.1573604060913706
$ws.Range("F17").Value = 0.01821862348178137
$ws.Range("H17").Value = 0.1740890688259109
$ws.Range("I17").Value = 0.08097165991902834
$ws.Range("J17").Value = 0.4676113360323887
$ws.Range("K17").Value = 0.07692307692307693
$ws.Range("M17").Value = 0.01821862348178137
$ws.Range("O17").Value = 0.05668016194331984
$ws.Range("S17").Value = 0.1072874493927125
$ws.Range("F18").Value = 0.01858736059479554
$ws.Range("H18").Value = 0.1858736059479554
$ws.Range("I18").Value = 0.1078066914498141
$ws.Range("J18").Value = 0.412639405204461
$ws.Range("K18").Value = 0.07063197026022305
$ws.Range("M18").Value = 0.02973977695167286
$ws.Range("O18").Value = 0.07063197026022305
$ws.Range("S18").Value = 0.104089219330855
$ws.Range("F19").Value = 0.0185334407735697
$ws.Range("H19").Value = 0.1990330378726833
$ws.Range("I19").Value = 0.0926672038678485
$ws.Range("J19").Value = 0.4045124899274778
$ws.Range("K19").Value = 0.08138597904915391
$ws.Range("M19").Value = 0.0290088638195004
$ws.Range("N19").Value = 0.0008058017727639
$ws.Range("O19").Value = 0.08058017727639001
$ws.Range("S19").Value = 0.09347300564061241
